$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before the existing "Big John" rows (currently rows 145-146),
# pushing them down to rows 150-151.
$ws.Rows("145:149").Insert()

# Common values shared across all the new rows (145-149)
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$fecha     = 44595
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100103
$producto  = "Frutos de hueso (carozo)"
$categoriaId = 100103006
$categoria = "Nectarín"
$unidad    = "`$/bins (420 kilos)"
$origen    = "Región de O'Higgins"
$kgUnidad  = 420

$rows = @(
    @{ Row = 145; Variedad = "June Pearl"; Calidad = "Especial"; Volumen = 20; Min = 390000; Max = 400000; Prom = 395000; PrecioKg = 940 },
    @{ Row = 146; Variedad = "June Pearl"; Calidad = "Primera";  Volumen = 20; Min = 350000; Max = 360000; Prom = 355000; PrecioKg = 845 },
    @{ Row = 147; Variedad = "Venus";      Calidad = "Especial"; Volumen = 16; Min = 360000; Max = 370000; Prom = 365000; PrecioKg = 869 },
    @{ Row = 148; Variedad = "Venus";      Calidad = "Primera";  Volumen = 20; Min = 320000; Max = 330000; Prom = 325000; PrecioKg = 774 },
    @{ Row = 149; Variedad = "Venus";      Calidad = "Segunda";  Volumen = 20; Min = 270000; Max = 280000; Prom = 275000; PrecioKg = 655 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.Min
    $ws.Cells.Item($row, 15).Value = $r.Max
    $ws.Cells.Item($row, 16).Value = $r.Prom
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
